$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.098.37'
$ws.Range('E2').Value = '  -1.56%  '
$ws.Range('D3').Value = '1.798.46'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '222.61'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.20'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.44%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.284'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0714'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.19%  '
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('D12').Value = '2.055.77'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '1.795.62'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.69'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.630'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.44%  '
$ws.Range('D16').Value = '34.137.44'
$ws.Range('E16').Value = '  -1.57%  '
$ws.Range('E17').Value = '  -1.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.13'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '246.48'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.48%  '
$ws.Range('D20').Value = '0.0₃0783'
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.83'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +4.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.08'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.11'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.89'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.55'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.05'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.33%  '
$ws.Range('E28').Value = '  -1.26%  '
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').Value = '  +1.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.71'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.47%  '
$ws.Range('E32').Value = '  +1.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.50'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.84'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.37%  '
$ws.Range('D35').Value = '1.411.99'
$ws.Range('E35').Value = '  -0.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.644'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.23%  '
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0186'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.942'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '80.20'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.26%  '
$ws.Range('E41').Value = '  -2.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.34'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.14'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +4.66%  '
$ws.Range('E44').Value = '  +0.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0496'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '107.19'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.54%  '
$ws.Range('D47').Value = '1.954.87'
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('E48').Value = '  -2.72%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.93'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.07%  '
$ws.Range('D51').Value = '0.0₆0123'
$ws.Range('E51').Value = '  +1.78%  '
